# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" sheet (sheet1) and the aggregate "全部类型" sheet (sheet4),
# reflecting freshly scraped numbers as of commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 3217
$ws1.Range("F7").Value = 748
$ws1.Range("F8").Value = 1963
$ws1.Range("F16").Value = 55
$ws1.Range("F18").Value = 474
$ws1.Range("F19").Value = 580
$ws1.Range("F21").Value = 10345
$ws1.Range("F22").Value = 9540
$ws1.Range("F23").Value = 827
$ws1.Range("F25").Value = 1795
$ws1.Range("F26").Value = 126

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 3217
$ws4.Range("F9").Value = 748
$ws4.Range("F10").Value = 1963
$ws4.Range("F19").Value = 55
$ws4.Range("F22").Value = 474
$ws4.Range("F23").Value = 580
$ws4.Range("F25").Value = 10345
$ws4.Range("F26").Value = 9540
$ws4.Range("F27").Value = 827
$ws4.Range("F29").Value = 1795
$ws4.Range("F32").Value = 126
